$d = $word.ActiveDocument

# Mark every inline picture's run as "do not spell/grammar check" (w:noProof),
# matching Word's normal behavior when a picture is (re)inserted/refreshed.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = -1
}

# Append a new blank paragraph followed by a paragraph of text at the end
# of the document body (right before the section break).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter()

$newTextPara = $d.Paragraphs.Last
$newTextPara.Range.InsertBefore("Criar solution na pasta e depois cria a pasta src e adiciona os projetos la todos, pode ter sql tbm fora, testes e ...")
